$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column K (dHeight), shifting existing columns K..O to L..P
$ws.Range("K1:K3").EntireColumn.Insert()

# New column K: dULCarrierMHz header + values duplicated from column J (dDLCarrierMHz)
$ws.Range("K1").Value = "dULCarrierMHz"
$ws.Range("K2").Value = $ws.Range("J2").Value2
$ws.Range("K3").Value = $ws.Range("J3").Value2

# Rename header of the (now shifted) NodeType column (M1) to BSCat
$ws.Range("M1").Value = "BSCat"

# Update selection to match target (M1)
$ws.Range("M1").Select()
